$d = $word.ActiveDocument
$count = 0
$r = $d.Content
$res = $r.Find.Execute("28+48=", $true, $false, $false, $false, $false, $true, 1, $false, "52-19=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 28+48=" }
$r = $d.Content
$res = $r.Find.Execute("59-18=", $true, $false, $false, $false, $false, $true, 1, $false, "63-13=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 59-18=" }
$r = $d.Content
$res = $r.Find.Execute("25+3=", $true, $false, $false, $false, $false, $true, 1, $false, "58+4=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 25+3=" }
$r = $d.Content
$res = $r.Find.Execute("40-1=", $true, $false, $false, $false, $false, $true, 1, $false, "56-8=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 40-1=" }
$r = $d.Content
$res = $r.Find.Execute("5+4=", $true, $false, $false, $false, $false, $true, 1, $false, "49-26=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 5+4=" }
$r = $d.Content
$res = $r.Find.Execute("2+39=", $true, $false, $false, $false, $false, $true, 1, $false, "3+73=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 2+39=" }
$r = $d.Content
$res = $r.Find.Execute("30+23=", $true, $false, $false, $false, $false, $true, 1, $false, "83-7=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 30+23=" }
$r = $d.Content
$res = $r.Find.Execute("84+6=", $true, $false, $false, $false, $false, $true, 1, $false, "80-27=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 84+6=" }
$r = $d.Content
$res = $r.Find.Execute("83-9=", $true, $false, $false, $false, $false, $true, 1, $false, "57-31=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 83-9=" }
$r = $d.Content
$res = $r.Find.Execute("46-32=", $true, $false, $false, $false, $false, $true, 1, $false, "69-21=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 46-32=" }
$r = $d.Content
$res = $r.Find.Execute("51-26=", $true, $false, $false, $false, $false, $true, 1, $false, "25-5=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 51-26=" }
$r = $d.Content
$res = $r.Find.Execute("39-12=", $true, $false, $false, $false, $false, $true, 1, $false, "75-8=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 39-12=" }
$r = $d.Content
$res = $r.Find.Execute("34+16=", $true, $false, $false, $false, $false, $true, 1, $false, "51-23=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 34+16=" }
$r = $d.Content
$res = $r.Find.Execute("6+0=", $true, $false, $false, $false, $false, $true, 1, $false, "83-42=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 6+0=" }
$r = $d.Content
$res = $r.Find.Execute("62+17=", $true, $false, $false, $false, $false, $true, 1, $false, "62+15=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 62+17=" }
$r = $d.Content
$res = $r.Find.Execute("70-64=", $true, $false, $false, $false, $false, $true, 1, $false, "27+38=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 70-64=" }
$r = $d.Content
$res = $r.Find.Execute("33-7=", $true, $false, $false, $false, $false, $true, 1, $false, "67-29=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 33-7=" }
$r = $d.Content
$res = $r.Find.Execute("40+36=", $true, $false, $false, $false, $false, $true, 1, $false, "46+53=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 40+36=" }
$r = $d.Content
$res = $r.Find.Execute("51-33=", $true, $false, $false, $false, $false, $true, 1, $false, "25+24=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 51-33=" }
$r = $d.Content
$res = $r.Find.Execute("14+79=", $true, $false, $false, $false, $false, $true, 1, $false, "24-7=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 14+79=" }
$r = $d.Content
$res = $r.Find.Execute("13+12=", $true, $false, $false, $false, $false, $true, 1, $false, "47+33=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 13+12=" }
$r = $d.Content
$res = $r.Find.Execute("56+34=", $true, $false, $false, $false, $false, $true, 1, $false, "50-2=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 56+34=" }
$r = $d.Content
$res = $r.Find.Execute("23+7=", $true, $false, $false, $false, $false, $true, 1, $false, "49-27=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 23+7=" }
$r = $d.Content
$res = $r.Find.Execute("34+48=", $true, $false, $false, $false, $false, $true, 1, $false, "36+34=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 34+48=" }
$r = $d.Content
$res = $r.Find.Execute("76-51=", $true, $false, $false, $false, $false, $true, 1, $false, "48-23=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 76-51=" }
$r = $d.Content
$res = $r.Find.Execute("39+10=", $true, $false, $false, $false, $false, $true, 1, $false, "93-88=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 39+10=" }
$r = $d.Content
$res = $r.Find.Execute("51+26=", $true, $false, $false, $false, $false, $true, 1, $false, "15+48=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 51+26=" }
$r = $d.Content
$res = $r.Find.Execute("78-36=", $true, $false, $false, $false, $false, $true, 1, $false, "93-52=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 78-36=" }
$r = $d.Content
$res = $r.Find.Execute("38+61=", $true, $false, $false, $false, $false, $true, 1, $false, "15+60=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 38+61=" }
$r = $d.Content
$res = $r.Find.Execute("90-19=", $true, $false, $false, $false, $false, $true, 1, $false, "56-4=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 90-19=" }
$r = $d.Content
$res = $r.Find.Execute("39-35=", $true, $false, $false, $false, $false, $true, 1, $false, "76-37=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 39-35=" }
$r = $d.Content
$res = $r.Find.Execute("52+15=", $true, $false, $false, $false, $false, $true, 1, $false, "68-38=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 52+15=" }
$r = $d.Content
$res = $r.Find.Execute("65-13=", $true, $false, $false, $false, $false, $true, 1, $false, "98-31=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 65-13=" }
$r = $d.Content
$res = $r.Find.Execute("79-56=", $true, $false, $false, $false, $false, $true, 1, $false, "91-61=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 79-56=" }
$r = $d.Content
$res = $r.Find.Execute("63-49=", $true, $false, $false, $false, $false, $true, 1, $false, "47+22=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 63-49=" }
$r = $d.Content
$res = $r.Find.Execute("3+88=", $true, $false, $false, $false, $false, $true, 1, $false, "9+9=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 3+88=" }
$r = $d.Content
$res = $r.Find.Execute("62-3=", $true, $false, $false, $false, $false, $true, 1, $false, "90-18=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 62-3=" }
$r = $d.Content
$res = $r.Find.Execute("19+71=", $true, $false, $false, $false, $false, $true, 1, $false, "98-3=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 19+71=" }
$r = $d.Content
$res = $r.Find.Execute("58+5=", $true, $false, $false, $false, $false, $true, 1, $false, "83-43=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 58+5=" }
$r = $d.Content
$res = $r.Find.Execute("62-42=", $true, $false, $false, $false, $false, $true, 1, $false, "44+11=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 62-42=" }
$r = $d.Content
$res = $r.Find.Execute("14-9=", $true, $false, $false, $false, $false, $true, 1, $false, "80-9=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 14-9=" }
$r = $d.Content
$res = $r.Find.Execute("13+83=", $true, $false, $false, $false, $false, $true, 1, $false, "76-42=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 13+83=" }
$r = $d.Content
$res = $r.Find.Execute("75-35=", $true, $false, $false, $false, $false, $true, 1, $false, "81+1=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 75-35=" }
$r = $d.Content
$res = $r.Find.Execute("74-11=", $true, $false, $false, $false, $false, $true, 1, $false, "77-13=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 74-11=" }
$r = $d.Content
$res = $r.Find.Execute("88-75=", $true, $false, $false, $false, $false, $true, 1, $false, "49-19=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 88-75=" }
$r = $d.Content
$res = $r.Find.Execute("28+6=", $true, $false, $false, $false, $false, $true, 1, $false, "61-29=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 28+6=" }
$r = $d.Content
$res = $r.Find.Execute("90-88=", $true, $false, $false, $false, $false, $true, 1, $false, "14+27=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 90-88=" }
$r = $d.Content
$res = $r.Find.Execute("25+57=", $true, $false, $false, $false, $false, $true, 1, $false, "59+32=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 25+57=" }
$r = $d.Content
$res = $r.Find.Execute("43-3=", $true, $false, $false, $false, $false, $true, 1, $false, "36-6=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 43-3=" }
$r = $d.Content
$res = $r.Find.Execute("43+16=", $true, $false, $false, $false, $false, $true, 1, $false, "16-12=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 43+16=" }
$r = $d.Content
$res = $r.Find.Execute("16+14=", $true, $false, $false, $false, $false, $true, 1, $false, "19+46=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 16+14=" }
$r = $d.Content
$res = $r.Find.Execute("91-33=", $true, $false, $false, $false, $false, $true, 1, $false, "63-24=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 91-33=" }
$r = $d.Content
$res = $r.Find.Execute("39-22=", $true, $false, $false, $false, $false, $true, 1, $false, "38+34=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 39-22=" }
$r = $d.Content
$res = $r.Find.Execute("11+49=", $true, $false, $false, $false, $false, $true, 1, $false, "43-39=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 11+49=" }
$r = $d.Content
$res = $r.Find.Execute("20-15=", $true, $false, $false, $false, $false, $true, 1, $false, "50+37=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 20-15=" }
$r = $d.Content
$res = $r.Find.Execute("10+26=", $true, $false, $false, $false, $false, $true, 1, $false, "48+33=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 10+26=" }
$r = $d.Content
$res = $r.Find.Execute("38+43=", $true, $false, $false, $false, $false, $true, 1, $false, "23+29=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 38+43=" }
$r = $d.Content
$res = $r.Find.Execute("85-30=", $true, $false, $false, $false, $false, $true, 1, $false, "59-47=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 85-30=" }
$r = $d.Content
$res = $r.Find.Execute("91-68=", $true, $false, $false, $false, $false, $true, 1, $false, "26-23=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 91-68=" }
$r = $d.Content
$res = $r.Find.Execute("32+0=", $true, $false, $false, $false, $false, $true, 1, $false, "99-51=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 32+0=" }
$r = $d.Content
$res = $r.Find.Execute("73+1=", $true, $false, $false, $false, $false, $true, 1, $false, "96-49=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 73+1=" }
$r = $d.Content
$res = $r.Find.Execute("8+81=", $true, $false, $false, $false, $false, $true, 1, $false, "2+62=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 8+81=" }
$r = $d.Content
$res = $r.Find.Execute("81-12=", $true, $false, $false, $false, $false, $true, 1, $false, "0+69=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 81-12=" }
$r = $d.Content
$res = $r.Find.Execute("87-14=", $true, $false, $false, $false, $false, $true, 1, $false, "90-43=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 87-14=" }
$r = $d.Content
$res = $r.Find.Execute("72+11=", $true, $false, $false, $false, $false, $true, 1, $false, "18+39=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 72+11=" }
$r = $d.Content
$res = $r.Find.Execute("35+9=", $true, $false, $false, $false, $false, $true, 1, $false, "55-46=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 35+9=" }
$r = $d.Content
$res = $r.Find.Execute("91-16=", $true, $false, $false, $false, $false, $true, 1, $false, "32-11=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 91-16=" }
$r = $d.Content
$res = $r.Find.Execute("77-62=", $true, $false, $false, $false, $false, $true, 1, $false, "53-2=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 77-62=" }
$r = $d.Content
$res = $r.Find.Execute("92-8=", $true, $false, $false, $false, $false, $true, 1, $false, "57+39=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 92-8=" }
$r = $d.Content
$res = $r.Find.Execute("51-13=", $true, $false, $false, $false, $false, $true, 1, $false, "55+2=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 51-13=" }
$r = $d.Content
$res = $r.Find.Execute("59+2=", $true, $false, $false, $false, $false, $true, 1, $false, "68-23=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 59+2=" }
$r = $d.Content
$res = $r.Find.Execute("39-27=", $true, $false, $false, $false, $false, $true, 1, $false, "18+70=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 39-27=" }
$r = $d.Content
$res = $r.Find.Execute("34+55=", $true, $false, $false, $false, $false, $true, 1, $false, "79-33=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 34+55=" }
$r = $d.Content
$res = $r.Find.Execute("32-2=", $true, $false, $false, $false, $false, $true, 1, $false, "75+9=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 32-2=" }
$r = $d.Content
$res = $r.Find.Execute("78-43=", $true, $false, $false, $false, $false, $true, 1, $false, "71+12=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 78-43=" }
$r = $d.Content
$res = $r.Find.Execute("71-24=", $true, $false, $false, $false, $false, $true, 1, $false, "51-20=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 71-24=" }
$r = $d.Content
$res = $r.Find.Execute("43-23=", $true, $false, $false, $false, $false, $true, 1, $false, "10+11=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 43-23=" }
$r = $d.Content
$res = $r.Find.Execute("93-17=", $true, $false, $false, $false, $false, $true, 1, $false, "82+17=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 93-17=" }
$r = $d.Content
$res = $r.Find.Execute("57-44=", $true, $false, $false, $false, $false, $true, 1, $false, "43+6=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 57-44=" }
$r = $d.Content
$res = $r.Find.Execute("81-17=", $true, $false, $false, $false, $false, $true, 1, $false, "70+2=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 81-17=" }
$r = $d.Content
$res = $r.Find.Execute("54-19=", $true, $false, $false, $false, $false, $true, 1, $false, "47+21=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 54-19=" }
$r = $d.Content
$res = $r.Find.Execute("51-14=", $true, $false, $false, $false, $false, $true, 1, $false, "69-60=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 51-14=" }
$r = $d.Content
$res = $r.Find.Execute("91-46=", $true, $false, $false, $false, $false, $true, 1, $false, "81-79=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 91-46=" }
$r = $d.Content
$res = $r.Find.Execute("55+11=", $true, $false, $false, $false, $false, $true, 1, $false, "93-72=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 55+11=" }
$r = $d.Content
$res = $r.Find.Execute("85-36=", $true, $false, $false, $false, $false, $true, 1, $false, "35+15=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 85-36=" }
$r = $d.Content
$res = $r.Find.Execute("64-9=", $true, $false, $false, $false, $false, $true, 1, $false, "47+52=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 64-9=" }
$r = $d.Content
$res = $r.Find.Execute("56-56=", $true, $false, $false, $false, $false, $true, 1, $false, "8+34=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 56-56=" }
$r = $d.Content
$res = $r.Find.Execute("11+42=", $true, $false, $false, $false, $false, $true, 1, $false, "46+2=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 11+42=" }
$r = $d.Content
$res = $r.Find.Execute("52+31=", $true, $false, $false, $false, $false, $true, 1, $false, "21-2=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 52+31=" }
$r = $d.Content
$res = $r.Find.Execute("51+46=", $true, $false, $false, $false, $false, $true, 1, $false, "48+49=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 51+46=" }
$r = $d.Content
$res = $r.Find.Execute("26+54=", $true, $false, $false, $false, $false, $true, 1, $false, "3+91=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 26+54=" }
$r = $d.Content
$res = $r.Find.Execute("23+69=", $true, $false, $false, $false, $false, $true, 1, $false, "23-20=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 23+69=" }
$r = $d.Content
$res = $r.Find.Execute("89-84=", $true, $false, $false, $false, $false, $true, 1, $false, "15+67=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 89-84=" }
$r = $d.Content
$res = $r.Find.Execute("58+14=", $true, $false, $false, $false, $false, $true, 1, $false, "20+1=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 58+14=" }
$r = $d.Content
$res = $r.Find.Execute("51-27=", $true, $false, $false, $false, $false, $true, 1, $false, "20+13=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 51-27=" }
$r = $d.Content
$res = $r.Find.Execute("37+21=", $true, $false, $false, $false, $false, $true, 1, $false, "48-19=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 37+21=" }
$r = $d.Content
$res = $r.Find.Execute("12+16=", $true, $false, $false, $false, $false, $true, 1, $false, "33+11=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 12+16=" }
$r = $d.Content
$res = $r.Find.Execute("87-22=", $true, $false, $false, $false, $false, $true, 1, $false, "73+6=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 87-22=" }
$r = $d.Content
$res = $r.Find.Execute("13+36=", $true, $false, $false, $false, $false, $true, 1, $false, "34+56=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 13+36=" }
$r = $d.Content
$res = $r.Find.Execute("67-38=", $true, $false, $false, $false, $false, $true, 1, $false, "76-14=", 2)
if ($res) { $count = $count + 1 } else { Write-Output "MISSING: 67-38=" }
Write-Output "Replaced $count of 100 "
